# Fix Training Data Issue
# The "Date" column (BF) on Sheet1 was recorded one day off because of how
# the NBA stats site displayed its dates (e.g. "4-7-2013-14" instead of the
# correct ISO style date "2014-04-07"). Correct every data row (rows 2-31)
# in column BF ("Date") so the model is trained on the right date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "4-7-2013-14"
$newDate = "2014-04-07"

$dateColumn = 58   # column BF
$firstDataRow = 2
$lastDataRow = 31

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateColumn)

    if ($cell.Value2 -eq $oldDate) {
        # Force the cell to stay text (instead of Excel auto-converting the
        # ISO-looking string into a date serial number) before writing the
        # corrected value.
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
    }
}
